$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 'Zagłębie Lubin'
$ws.Range("C11").Value = 'Remis'
$ws.Range("C15").Value = 'Zagłębie Lubin'
$ws.Range("C17").Value = 'Radomiak Radom'
$ws.Range("C18").Value = 'Remis'
$ws.Range("C20").Value = 'Legia Warszawa'
$ws.Range("C26").Value = 'Górnik Zabrze'
$ws.Range("C27").Value = 'Remis'
$ws.Range("C29").Value = 'Jagielonia Białystok'
$ws.Range("C32").Value = 'Legia Warszawa'
$ws.Range("C35").Value = 'Remis'
$ws.Range("C39").Value = 'Warta Poznań'
$ws.Range("C40").Value = 'Remis'
$ws.Range("C41").Value = 'Miedź Legnica'
$ws.Range("C42").Value = 'Górnik Zabrze'
$ws.Range("C47").Value = 'Stal Mielec'
$ws.Range("C48").Value = 'Jagielonia Białystok'
$ws.Range("C58").Value = 'Lechia Gdańsk'
$ws.Range("C60").Value = 'Pogoń Szczecin'
$ws.Range("C62").Value = 'Legia Warszawa'
$ws.Range("C63").Value = 'Remis'
$ws.Range("C65").Value = 'Remis'
$ws.Range("C66").Value = 'Miedź Legnica'
$ws.Range("C68").Value = 'Remis'
$ws.Range("C72").Value = 'Remis'
$ws.Range("C74").Value = 'Jagielonia Białystok'
$ws.Range("C75").Value = 'Remis'
$ws.Range("C78").Value = 'Remis'
$ws.Range("C79").Value = 'Remis'
$ws.Range("C80").Value = 'Remis'
$ws.Range("C81").Value = 'Remis'
$ws.Range("C82").Value = 'Śląsk Wrocław'
$ws.Range("C84").Value = 'Piast Gliwice'
$ws.Range("C88").Value = 'Miedź Legnica'
$ws.Range("C91").Value = 'Lech Poznań'
$ws.Range("C92").Value = 'Korona Kielce'
$ws.Range("C93").Value = 'Legia Warszawa'
$ws.Range("C95").Value = 'Remis'
$ws.Range("C96").Value = 'Remis'
$ws.Range("C99").Value = 'Wisła Płock'
$ws.Range("C100").Value = 'Remis'
$ws.Range("C101").Value = 'Remis'
$ws.Range("C103").Value = 'Remis'
$ws.Range("C108").Value = 'Pogoń Szczecin'
$ws.Range("C109").Value = 'Remis'
$ws.Range("C111").Value = 'Remis'
$ws.Range("C112").Value = 'Cracovia'
$ws.Range("C113").Value = 'Remis'
$ws.Range("C115").Value = 'Śląsk Wrocław'
$ws.Range("C117").Value = 'Remis'
$ws.Range("C118").Value = 'Remis'
$ws.Range("C123").Value = 'Raków Częstochowa'
$ws.Range("C125").Value = 'Remis'
$ws.Range("C127").Value = 'Remis'
$ws.Range("C129").Value = 'Remis'
$ws.Range("C135").Value = 'Remis'
$ws.Range("C137").Value = 'Jagielonia Białystok'
$ws.Range("C142").Value = 'Pogoń Szczecin'
$ws.Range("C148").Value = 'Widzew Łódź'
$ws.Range("C154").Value = 'Legia Warszawa'
$ws.Range("C155").Value = 'Remis'
$ws.Range("C156").Value = 'Jagielonia Białystok'
$ws.Range("C159").Value = 'Radomiak Radom'
$ws.Range("C161").Value = 'Remis'
$ws.Range("C163").Value = 'Remis'
$ws.Range("C164").Value = 'Remis'
$ws.Range("C165").Value = 'Korona Kielce'
$ws.Range("C167").Value = 'Remis'
$ws.Range("C169").Value = 'Remis'
$ws.Range("C170").Value = 'Raków Częstochowa'
$ws.Range("C171").Value = 'Widzew Łódź'
$ws.Range("C174").Value = 'Remis'
$ws.Range("C177").Value = 'Radomiak Radom'
$ws.Range("C181").Value = 'Śląsk Wrocław'
$ws.Range("C182").Value = 'Cracovia'
$ws.Range("C184").Value = 'Lechia Gdańsk'
$ws.Range("C186").Value = 'Miedź Legnica'
$ws.Range("C189").Value = 'Remis'
$ws.Range("C190").Value = 'Śląsk Wrocław'
$ws.Range("C193").Value = 'Zagłębie Lubin'
$ws.Range("C196").Value = 'Stal Mielec'
$ws.Range("C200").Value = 'Remis'
$ws.Range("C201").Value = 'Remis'
$ws.Range("C202").Value = 'Lech Poznań'
$ws.Range("C203").Value = 'Miedź Legnica'
$ws.Range("C204").Value = 'Legia Warszawa'
$ws.Range("C206").Value = 'Remis'
$ws.Range("C208").Value = 'Warta Poznań'
$ws.Range("C210").Value = 'Górnik Zabrze'
$ws.Range("C211").Value = 'Pogoń Szczecin'
$ws.Range("C212").Value = 'Korona Kielce'
$ws.Range("C216").Value = 'Cracovia'
$ws.Range("C219").Value = 'Miedź Legnica'
$ws.Range("C222").Value = 'Legia Warszawa'
$ws.Range("C224").Value = 'Lechia Gdańsk'
$ws.Range("C227").Value = 'Widzew Łódź'
$ws.Range("C228").Value = 'Górnik Zabrze'
$ws.Range("C229").Value = 'Warta Poznań'
$ws.Range("C230").Value = 'Korona Kielce'
$ws.Range("C231").Value = 'Remis'
$ws.Range("C233").Value = 'Raków Częstochowa'
$ws.Range("C237").Value = 'Warta Poznań'
$ws.Range("C240").Value = 'Cracovia'
$ws.Range("C244").Value = 'Piast Gliwice'
$ws.Range("C245").Value = 'Remis'
$ws.Range("C246").Value = 'Piast Gliwice'
$ws.Range("C247").Value = 'Remis'
$ws.Range("C249").Value = 'Pogoń Szczecin'
$ws.Range("C250").Value = 'Lech Poznań'
$ws.Range("C252").Value = 'Remis'
$ws.Range("C254").Value = 'Remis'
$ws.Range("C255").Value = 'Lechia Gdańsk'
$ws.Range("C257").Value = 'Górnik Zabrze'
$ws.Range("C260").Value = 'Remis'
$ws.Range("C262").Value = 'Wisła Płock'
$ws.Range("C263").Value = 'Miedź Legnica'
$ws.Range("C265").Value = 'Jagielonia Białystok'
$ws.Range("C268").Value = 'Remis'
$ws.Range("C271").Value = 'Radomiak Radom'
$ws.Range("C275").Value = 'Zagłębie Lubin'
$ws.Range("C277").Value = 'Górnik Zabrze'
$ws.Range("C278").Value = 'Pogoń Szczecin'
$ws.Range("C279").Value = 'Remis'
$ws.Range("C280").Value = 'Wisła Płock'
$ws.Range("C282").Value = 'Remis'
$ws.Range("C283").Value = 'Jagielonia Białystok'
$ws.Range("C285").Value = 'Remis'
$ws.Range("C288").Value = 'Remis'
$ws.Range("C289").Value = 'Remis'
$ws.Range("C290").Value = 'Cracovia'
$ws.Range("C291").Value = 'Stal Mielec'
$ws.Range("C293").Value = 'Lechia Gdańsk'
$ws.Range("C294").Value = 'Pogoń Szczecin'
$ws.Range("C295").Value = 'Radomiak Radom'
$ws.Range("C299").Value = 'Remis'
$ws.Range("C302").Value = 'Remis'
$ws.Range("C303").Value = 'Miedź Legnica'
$ws.Range("C304").Value = 'Remis'
$ws.Range("C306").Value = 'Warta Poznań'
$ws.Range("C307").Value = 'Remis'
